# Regenerate save_data to use K (strikeouts) instead of Strike# for column G.
# Rows 2-20 on Sheet1 correspond to the 19 most recent starts; only the "K"
# column (G) values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 1
    4  = 2
    5  = 0
    6  = 3
    7  = 3
    8  = 5
    9  = 3
    10 = 3
    11 = 0
    12 = 5
    13 = 0
    14 = 2
    15 = 0
    16 = 1
    17 = 5
    18 = 5
    19 = 6
    20 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
